# Adapt column header formatting to respective input file names.
# - Rename "_old" header suffixes to "_FV2404"
# - Rename "_new" header suffixes to "_FV2410"
# - Turn the data range into a real Excel Table (ListObject)
# - Freeze the header row (pane split) on the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -----------------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns A..J hold the "_old" headers, columns L..U hold the "_new" headers.
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# --- 2. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into an Excel Table (ListObject) -----------
$tableRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

$wb.Save()
